# Apply the "exemplo xlsx rstac; ajuste faq" edits.
#
# Summary of changes:
#  - Informações Pertinentes!A2: "Temperatura" -> "Temperatura no banho de ar"
#  - Informações Pertinentes!B3: 54.5 -> 54.4
#  - Rastreabilidade!D2: "PR 465" -> "Lampe PR 465"
#  - Método de Medição!B1: header "equation" removed (cell cleared, style kept)
#  - Per-sheet active-cell selections updated (navigation state)
#  - Active sheet changes from "Metadados Principais" to "Software"

$wb = $excel.ActiveWorkbook

# --- Metadados Principais ---------------------------------------------
$ws = $wb.Worksheets.Item("Metadados Principais")
$ws.Range("B14").Select()

# --- Cliente -------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cliente")
$ws.Range("B3").Select()

# --- Informações Pertinentes ----------------------------------------------
$ws = $wb.Worksheets.Item("Informações Pertinentes")
$ws.Range("A2").Value = "Temperatura no banho de ar"
$ws.Range("B3").Value = 54.4
$ws.Range("E8").Select()

# --- Declarações -----------------------------------------------------------
$ws = $wb.Worksheets.Item("Declarações")
$ws.Range("B3").Select()

# --- Rastreabilidade -------------------------------------------------------
$ws = $wb.Worksheets.Item("Rastreabilidade")
$ws.Range("D2").Value = "Lampe PR 465"
$ws.Range("C12").Select()

# --- Método de Medição -------------------------------------------------------
$ws = $wb.Worksheets.Item("Método de Medição")
$ws.Range("B1").ClearContents()
$ws.Range("A17").Select()

# --- Mensurando --------------------------------------------------------------
$ws = $wb.Worksheets.Item("Mensurando")
$ws.Range("E3").Select()

# --- Índices -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Índices")
$ws.Range("A3").Select()

# --- Resultados ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Resultados")
$ws.Range("E3").Select()

# --- Observações -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Observações")
$ws.Range("D7").Select()

# --- Software (selected/activated last -> becomes the active tab) ----------------------
$ws = $wb.Worksheets.Item("Software")
$ws.Range("D16").Select()
